$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.286.57'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '3.054.40'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.36'
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.79'
$ws.Range("E6").Value = '  +3.94%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.046.85'
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.482'
$ws.Range("E12").Value = '  +5.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  -3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.41'
$ws.Range("E14").Value = '  +7.79%  '
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '66.307.86'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '3.547.92'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.24'
$ws.Range("E18").Value = '  +4.25%  '
$ws.Range("D19").Value = '3.042.96'
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.25'
$ws.Range("E20").Value = '  +16.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '469.74'
$ws.Range("E21").Value = '  +2.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.709'
$ws.Range("E22").Value = '  +2.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("E23").Value = '  +2.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.55'
$ws.Range("E24").Value = '  +1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.98'
$ws.Range("E25").Value = '  +5.10%  '
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.19'
$ws.Range("E27").Value = '  -3.45%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.41'
$ws.Range("E29").Value = '  +3.39%  '
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.64'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.119'
$ws.Range("E32").Value = '  +7.13%  '
$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0000101'
$ws.Range("E33").Value = '  -4.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.38'
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.87'
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.05'
$ws.Range("E38").Value = '  +9.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.06'
$ws.Range("E39").Value = '  -6.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.62'
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.311'
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.87'
$ws.Range("E42").Value = '  -5.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.121'
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.65'
$ws.Range("E44").Value = '  +2.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0361'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '385.86'
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("D47").Value = '2.746.82'
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.36'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.16'
$ws.Range("E50").Value = '  +5.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("E51").Value = '  +4.37%  '
